$d = $word.ActiveDocument

# Locate the "SKILLS" heading paragraph with Find (as suggested by the
# runtime's example), then walk forward two paragraphs: the first is the
# blank "Body" paragraph directly under the heading, the second is the
# first existing skill entry ("Python"). We insert the new "Bash" skill
# paragraph right before that "Python" paragraph, i.e. immediately after
# the blank paragraph - matching the diff.
$found = $d.Content
$ok = $found.Find.Execute("SKILLS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$skillsPara = $null
foreach ($p in $d.Paragraphs) {
    if (($p.Range.Start -le $found.Start) -and ($found.Start -lt $p.Range.End)) {
        $skillsPara = $p
        break
    }
}

$blankPara = $skillsPara.Next()
$firstSkillPara = $blankPara.Next()

$insertPoint = $d.Range($firstSkillPara.Range.Start, $firstSkillPara.Range.Start)

# Build the exact OOXML for the new paragraph (same Body style / run
# formatting as the other skill entries) and insert it via InsertXML so
# the paragraph mark + run properties match the author's edit precisely.
$bashParagraphXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr>' +
                '<w:pStyle w:val="Body"/>' +
                '<w:rPr>' +
                  '<w:shd w:val="nil" w:color="auto" w:fill="auto"/>' +
                  '<w:lang w:val="en-US"/>' +
                '</w:rPr>' +
              '</w:pPr>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:shd w:val="nil" w:color="auto" w:fill="auto"/>' +
                  '<w:rtl w:val="0"/>' +
                  '<w:lang w:val="en-US"/>' +
                '</w:rPr>' +
                '<w:t>Bash</w:t>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$insertPoint.InsertXML($bashParagraphXml)
